$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Shorten the "Naheffingsaanslag parkeerbelasting" bullet: drop the
#    trailing explanatory sentence, keep the lead-in clause intact.
# ---------------------------------------------------------------------------
$newText = "- Naheffingsaanslag parkeerbelasting (parkeerboete) - uitstel van betaling tot bezwaar is afgehandeld"

$naheffingParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Naheffingsaanslag parkeerbelasting*") {
        $naheffingParagraph = $p
        break
    }
}

if ($naheffingParagraph -ne $null) {
    $naheffingParagraph.Range.Text = $newText
}

# ---------------------------------------------------------------------------
# 2) Extend the WOZ bullet with ", WOZ waarde opvragen" as two additional
#    runs (matching the sz/szCs 26 formatting already used on the paragraph).
# ---------------------------------------------------------------------------
$wozParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*WOZ Inlichtingenformulier, WOZ taxatieverslag, WOZ-beschikking*") {
        $wozParagraph = $p
        break
    }
}

if ($wozParagraph -ne $null) {
    $paraRange = $wozParagraph.Range
    $insertPos = $paraRange.End - 1

    $run1 = $d.Range($insertPos, $insertPos)
    $run1.InsertAfter(", ")
    $run1.Font.Size = 13

    $insertPos2 = $insertPos + 2
    $run2 = $d.Range($insertPos2, $insertPos2)
    $run2.InsertAfter("WOZ waarde opvragen")
    $run2.Font.Size = 13
}
